# Applies the diff: data/cleaned_data/daily_profile_metrics.xlsx
# - corrects the precision of G33 (extraction_datetime)
# - appends 10 new daily rows (34-43) with the next scrape results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: extraction_datetime precision correction ---
$ws.Cells.Item(33, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(33, 7).Value = 45779.39787231482

# --- Row 34 ---
$ws.Cells.Item(34, 1).Value = "ʜᴇʟᴘɪɴɢ ʏᴏᴜ ᴍᴏᴠᴇ ғʀᴏᴍ ʜᴇsɪᴛᴀᴛɪᴏɴ ᴛᴏ ᴄʀᴇᴀᴛɪᴏɴ
✨ | ✧ 𝗹𝗼𝗰𝘀 ✧ (𝘀𝗲𝗹𝗳) 𝗹𝗶𝗯𝗲𝗿𝗮𝘁𝗶𝗼𝗻 ✧ 𝗹𝗶𝗳𝗲𝘀𝘁𝘆𝗹𝗲 ✧
🪴 | 71 ʟᴏᴄs est. on 07.20.23
📍 | ʜᴏᴜsᴛᴏɴ, ᴛx"
$ws.Cells.Item(34, 2).Value = 2947
$ws.Cells.Item(34, 3).Value = 237
$ws.Cells.Item(34, 4).Value = 166
$ws.Cells.Item(34, 5).Value = "https://scontent-hou1-1.xx.fbcdn.net/v/t51.2885-15/481266977_997353345602937_1719041919639027270_n.jpg?_nc_cat=106&ccb=1-7&_nc_sid=7d201b&_nc_ohc=IN_GzJ2sEioQ7kNvwEqMVvT&_nc_oc=Adn9AyAP7zVKP8mCM-NlEfWXWDn0EUPNIDmmlWymEbMxy_VQTH1j2nzPajVTTI-IyX7LCUt3-J_BFTygH72RmzK3&_nc_zt=23&_nc_ht=scontent-hou1-1.xx&edm=AL-3X8kEAAAA&oh=00_AfHdpU4Tt9hgkmkwxmNbcV03t3tVxjWexSSwb3ifD4AmOQ&oe=681BF46E"
$ws.Cells.Item(34, 6).Value = "'17841461458191255"
$ws.Cells.Item(34, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(34, 7).Value = 45780.39800649306
$ws.Cells.Item(34, 8).Value = "'2025-05-03"
$ws.Cells.Item(34, 9).Value = 2025
$ws.Cells.Item(34, 10).Value = "May"
$ws.Cells.Item(34, 11).Value = 3
$ws.Cells.Item(34, 12).Value = "09:33:07"

# --- Row 35 ---
$ws.Cells.Item(35, 1).Value = "ʜᴇʟᴘɪɴɢ ʏᴏᴜ ᴍᴏᴠᴇ ғʀᴏᴍ ʜᴇsɪᴛᴀᴛɪᴏɴ ᴛᴏ ᴄʀᴇᴀᴛɪᴏɴ
✨ | ✧ 𝗹𝗼𝗰𝘀 ✧ (𝘀𝗲𝗹𝗳) 𝗹𝗶𝗯𝗲𝗿𝗮𝘁𝗶𝗼𝗻 ✧ 𝗹𝗶𝗳𝗲𝘀𝘁𝘆𝗹𝗲 ✧
🪴 | 71 ʟᴏᴄs est. on 07.20.23
📍 | ʜᴏᴜsᴛᴏɴ, ᴛx"
$ws.Cells.Item(35, 2).Value = 2950
$ws.Cells.Item(35, 3).Value = 237
$ws.Cells.Item(35, 4).Value = 167
$ws.Cells.Item(35, 5).Value = "https://scontent-hou1-1.xx.fbcdn.net/v/t51.2885-15/481266977_997353345602937_1719041919639027270_n.jpg?_nc_cat=106&ccb=1-7&_nc_sid=7d201b&_nc_ohc=IN_GzJ2sEioQ7kNvwEg5sAD&_nc_oc=AdmbXpAWtsQMCdIwbFTEUc4JIrzUeJ9uIJ4C1O82AyJcxE37ZsWHolsvcoPOTNlvUekXXV9pY0xppZrwjtvdwQuG&_nc_zt=23&_nc_ht=scontent-hou1-1.xx&edm=AL-3X8kEAAAA&oh=00_AfFQOke_4Tg7xiVQUNcuQXo67qvYOZgww0dl-EkFtkq7Ww&oe=681D45EE"
$ws.Cells.Item(35, 6).Value = "'17841461458191255"
$ws.Cells.Item(35, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(35, 7).Value = 45781.39785303241
$ws.Cells.Item(35, 8).Value = "'2025-05-04"
$ws.Cells.Item(35, 9).Value = 2025
$ws.Cells.Item(35, 10).Value = "May"
$ws.Cells.Item(35, 11).Value = 4
$ws.Cells.Item(35, 12).Value = "09:32:54"

# --- Row 36 ---
$ws.Cells.Item(36, 1).Value = "ʜᴇʟᴘɪɴɢ ʏᴏᴜ ᴍᴏᴠᴇ ғʀᴏᴍ ʜᴇsɪᴛᴀᴛɪᴏɴ ᴛᴏ ᴄʀᴇᴀᴛɪᴏɴ
✨ | ✧ 𝗹𝗼𝗰𝘀 ✧ (𝘀𝗲𝗹𝗳) 𝗹𝗶𝗯𝗲𝗿𝗮𝘁𝗶𝗼𝗻 ✧ 𝗹𝗶𝗳𝗲𝘀𝘁𝘆𝗹𝗲 ✧
🪴 | 71 ʟᴏᴄs est. on 07.20.23
📍 | ʜᴏᴜsᴛᴏɴ, ᴛx"
$ws.Cells.Item(36, 2).Value = 2952
$ws.Cells.Item(36, 3).Value = 237
$ws.Cells.Item(36, 4).Value = 167
$ws.Cells.Item(36, 5).Value = "https://scontent-hou1-1.xx.fbcdn.net/v/t51.2885-15/481266977_997353345602937_1719041919639027270_n.jpg?_nc_cat=106&ccb=1-7&_nc_sid=7d201b&_nc_ohc=NWtNeasGJ7MQ7kNvwFv6SEU&_nc_oc=Adl4mAinciDaL3ehOoV15zVu7rknF2pd7VX55SvNk7sVYaKBpUDwGi7Jq52lF0FmLp7QfmpP3mLsnvLUT7T1A-Sm&_nc_zt=23&_nc_ht=scontent-hou1-1.xx&edm=AL-3X8kEAAAA&oh=00_AfHcbyEFnu0D8TK_-SHc-MnU6EbB3y2P69duW_FBMW99bA&oe=681E976E"
$ws.Cells.Item(36, 6).Value = "'17841461458191255"
$ws.Cells.Item(36, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(36, 7).Value = 45782.40760672454
$ws.Cells.Item(36, 8).Value = "'2025-05-05"
$ws.Cells.Item(36, 9).Value = 2025
$ws.Cells.Item(36, 10).Value = "May"
$ws.Cells.Item(36, 11).Value = 5
$ws.Cells.Item(36, 12).Value = "09:46:57"

# --- Row 37 ---
$ws.Cells.Item(37, 1).Value = "ʜᴇʟᴘɪɴɢ ʏᴏᴜ ᴍᴏᴠᴇ ғʀᴏᴍ ʜᴇsɪᴛᴀᴛɪᴏɴ ᴛᴏ ᴄʀᴇᴀᴛɪᴏɴ
✨ | ✧ 𝗹𝗼𝗰𝘀 ✧ (𝘀𝗲𝗹𝗳) 𝗹𝗶𝗯𝗲𝗿𝗮𝘁𝗶𝗼𝗻 ✧ 𝗹𝗶𝗳𝗲𝘀𝘁𝘆𝗹𝗲 ✧
🪴 | 71 ʟᴏᴄs est. on 07.20.23
📍 | ʜᴏᴜsᴛᴏɴ, ᴛx"
$ws.Cells.Item(37, 2).Value = 2954
$ws.Cells.Item(37, 3).Value = 237
$ws.Cells.Item(37, 4).Value = 167
$ws.Cells.Item(37, 5).Value = "https://scontent-hou1-1.xx.fbcdn.net/v/t51.2885-15/481266977_997353345602937_1719041919639027270_n.jpg?_nc_cat=106&ccb=1-7&_nc_sid=7d201b&_nc_ohc=NWtNeasGJ7MQ7kNvwFv6SEU&_nc_oc=Adl4mAinciDaL3ehOoV15zVu7rknF2pd7VX55SvNk7sVYaKBpUDwGi7Jq52lF0FmLp7QfmpP3mLsnvLUT7T1A-Sm&_nc_zt=23&_nc_ht=scontent-hou1-1.xx&edm=AL-3X8kEAAAA&oh=00_AfKR6iib3jExJe38KR9WHOoXEmZ7y1ZaJmiaEMtZ4Q4gYQ&oe=681FE8EE"
$ws.Cells.Item(37, 6).Value = "'17841461458191255"
$ws.Cells.Item(37, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(37, 7).Value = 45783.39771293981
$ws.Cells.Item(37, 8).Value = "'2025-05-06"
$ws.Cells.Item(37, 9).Value = 2025
$ws.Cells.Item(37, 10).Value = "May"
$ws.Cells.Item(37, 11).Value = 6
$ws.Cells.Item(37, 12).Value = "09:32:42"

# --- Row 38 ---
$ws.Cells.Item(38, 1).Value = "ʜᴇʟᴘɪɴɢ ʏᴏᴜ ᴍᴏᴠᴇ ғʀᴏᴍ ʜᴇsɪᴛᴀᴛɪᴏɴ ᴛᴏ ᴄʀᴇᴀᴛɪᴏɴ
✨ | ✧ 𝗹𝗼𝗰𝘀 ✧ (𝘀𝗲𝗹𝗳) 𝗹𝗶𝗯𝗲𝗿𝗮𝘁𝗶𝗼𝗻 ✧ 𝗹𝗶𝗳𝗲𝘀𝘁𝘆𝗹𝗲 ✧
🪴 | 71 ʟᴏᴄs est. on 07.20.23
📍 | ʜᴏᴜsᴛᴏɴ, ᴛx"
$ws.Cells.Item(38, 2).Value = 2962
$ws.Cells.Item(38, 3).Value = 237
$ws.Cells.Item(38, 4).Value = 168
$ws.Cells.Item(38, 5).Value = "https://scontent-hou1-1.xx.fbcdn.net/v/t51.2885-15/481266977_997353345602937_1719041919639027270_n.jpg?_nc_cat=106&ccb=1-7&_nc_sid=7d201b&_nc_ohc=NWtNeasGJ7MQ7kNvwFv6SEU&_nc_oc=Adl4mAinciDaL3ehOoV15zVu7rknF2pd7VX55SvNk7sVYaKBpUDwGi7Jq52lF0FmLp7QfmpP3mLsnvLUT7T1A-Sm&_nc_zt=23&_nc_ht=scontent-hou1-1.xx&edm=AL-3X8kEAAAA&oh=00_AfKhIWqa2E3vK2KuY1WvA1T7IkZMBZDSDA2J_HgJZp_hXw&oe=68213A6E"
$ws.Cells.Item(38, 6).Value = "'17841461458191255"
$ws.Cells.Item(38, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(38, 7).Value = 45784.39764586806
$ws.Cells.Item(38, 8).Value = "'2025-05-07"
$ws.Cells.Item(38, 9).Value = 2025
$ws.Cells.Item(38, 10).Value = "May"
$ws.Cells.Item(38, 11).Value = 7
$ws.Cells.Item(38, 12).Value = "09:32:36"

# --- Row 39 ---
$ws.Cells.Item(39, 1).Value = "ʜᴇʟᴘɪɴɢ ʏᴏᴜ ᴍᴏᴠᴇ ғʀᴏᴍ ʜᴇsɪᴛᴀᴛɪᴏɴ ᴛᴏ ᴄʀᴇᴀᴛɪᴏɴ
✨ | ✧ 𝗹𝗼𝗰𝘀 ✧ (𝘀𝗲𝗹𝗳) 𝗹𝗶𝗯𝗲𝗿𝗮𝘁𝗶𝗼𝗻 ✧ 𝗹𝗶𝗳𝗲𝘀𝘁𝘆𝗹𝗲 ✧
🪴 | 71 ʟᴏᴄs est. on 07.20.23
📍 | ʜᴏᴜsᴛᴏɴ, ᴛx"
$ws.Cells.Item(39, 2).Value = 2966
$ws.Cells.Item(39, 3).Value = 237
$ws.Cells.Item(39, 4).Value = 168
$ws.Cells.Item(39, 5).Value = "https://scontent-hou1-1.xx.fbcdn.net/v/t51.2885-15/481266977_997353345602937_1719041919639027270_n.jpg?_nc_cat=106&ccb=1-7&_nc_sid=7d201b&_nc_ohc=0AU5y29RBJMQ7kNvwG0fiXo&_nc_oc=Adm6mtt-VWMmU47z_8dZcs8ZohztvDlKlvxT8c_mrNeBEpKue1MFDM6mNXCVGzOu5cmSkClZFSB2bAzP98sevVfi&_nc_zt=23&_nc_ht=scontent-hou1-1.xx&edm=AL-3X8kEAAAA&oh=00_AfJhxJX46NBaVXEX2MrzjBlSB1VdmWaCba61JaW8dMi86A&oe=68228BEE"
$ws.Cells.Item(39, 6).Value = "'17841461458191255"
$ws.Cells.Item(39, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(39, 7).Value = 45785.39818737269
$ws.Cells.Item(39, 8).Value = "'2025-05-08"
$ws.Cells.Item(39, 9).Value = 2025
$ws.Cells.Item(39, 10).Value = "May"
$ws.Cells.Item(39, 11).Value = 8
$ws.Cells.Item(39, 12).Value = "09:33:23"

# --- Row 40 ---
$ws.Cells.Item(40, 1).Value = "ʜᴇʟᴘɪɴɢ ʏᴏᴜ ᴍᴏᴠᴇ ғʀᴏᴍ ʜᴇsɪᴛᴀᴛɪᴏɴ ᴛᴏ ᴄʀᴇᴀᴛɪᴏɴ
✨ | ✧ 𝗹𝗼𝗰𝘀 ✧ (𝘀𝗲𝗹𝗳) 𝗹𝗶𝗯𝗲𝗿𝗮𝘁𝗶𝗼𝗻 ✧ 𝗹𝗶𝗳𝗲𝘀𝘁𝘆𝗹𝗲 ✧
🪴 | 71 ʟᴏᴄs est. on 07.20.23
📍 | ʜᴏᴜsᴛᴏɴ, ᴛx"
$ws.Cells.Item(40, 2).Value = 2970
$ws.Cells.Item(40, 3).Value = 238
$ws.Cells.Item(40, 4).Value = 169
$ws.Cells.Item(40, 5).Value = "https://scontent-hou1-1.xx.fbcdn.net/v/t51.2885-15/481266977_997353345602937_1719041919639027270_n.jpg?_nc_cat=106&ccb=1-7&_nc_sid=7d201b&_nc_ohc=0AU5y29RBJMQ7kNvwFl-fOI&_nc_oc=AdmC-OprH3Cbhp3s-2CBBrKKIWxNfyNh-7gYSoCwds5tNlOSVtv23jaLv7QJoH9Oqp-fKP8TUYosaQdDqdrVHnqj&_nc_zt=23&_nc_ht=scontent-hou1-1.xx&edm=AL-3X8kEAAAA&oh=00_AfJWCt6ApHgyf78T5VAwlD277JTJgC_WJVb--fiocoFkxA&oe=6823DD6E"
$ws.Cells.Item(40, 6).Value = "'17841461458191255"
$ws.Cells.Item(40, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(40, 7).Value = 45786.39813979167
$ws.Cells.Item(40, 8).Value = "'2025-05-09"
$ws.Cells.Item(40, 9).Value = 2025
$ws.Cells.Item(40, 10).Value = "May"
$ws.Cells.Item(40, 11).Value = 9
$ws.Cells.Item(40, 12).Value = "09:33:19"

# --- Row 41 ---
$ws.Cells.Item(41, 1).Value = "ʜᴇʟᴘɪɴɢ ʏᴏᴜ ᴍᴏᴠᴇ ғʀᴏᴍ ʜᴇsɪᴛᴀᴛɪᴏɴ ᴛᴏ ᴄʀᴇᴀᴛɪᴏɴ
✨ | ✧ 𝗹𝗼𝗰𝘀 ✧ (𝘀𝗲𝗹𝗳) 𝗹𝗶𝗯𝗲𝗿𝗮𝘁𝗶𝗼𝗻 ✧ 𝗹𝗶𝗳𝗲𝘀𝘁𝘆𝗹𝗲 ✧
🪴 | 71 ʟᴏᴄs est. on 07.20.23
📍 | ʜᴏᴜsᴛᴏɴ, ᴛx"
$ws.Cells.Item(41, 2).Value = 2973
$ws.Cells.Item(41, 3).Value = 238
$ws.Cells.Item(41, 4).Value = 170
$ws.Cells.Item(41, 5).Value = "https://scontent-hou1-1.xx.fbcdn.net/v/t51.2885-15/481266977_997353345602937_1719041919639027270_n.jpg?_nc_cat=106&ccb=1-7&_nc_sid=7d201b&_nc_ohc=0AU5y29RBJMQ7kNvwFpX11s&_nc_oc=AdlpXoSD_0bybCq1WNU1em450t9VxY0rAF519qzRUMnNfHN8bjVZtp-epfLj91LZ1Ub5tE3yQZi2vUG5XqwD9MCS&_nc_zt=23&_nc_ht=scontent-hou1-1.xx&edm=AL-3X8kEAAAA&oh=00_AfIRWnTm4YydRQ4X3NssjvKZfbW9KFgSkGxIxdremZ9JZg&oe=68252EEE"
$ws.Cells.Item(41, 6).Value = "'17841461458191255"
$ws.Cells.Item(41, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(41, 7).Value = 45787.39810366898
$ws.Cells.Item(41, 8).Value = "'2025-05-10"
$ws.Cells.Item(41, 9).Value = 2025
$ws.Cells.Item(41, 10).Value = "May"
$ws.Cells.Item(41, 11).Value = 10
$ws.Cells.Item(41, 12).Value = "09:33:16"

# --- Row 42 ---
$ws.Cells.Item(42, 1).Value = "ʜᴇʟᴘɪɴɢ ʏᴏᴜ ᴍᴏᴠᴇ ғʀᴏᴍ ʜᴇsɪᴛᴀᴛɪᴏɴ ᴛᴏ ᴄʀᴇᴀᴛɪᴏɴ
✨ | ✧ 𝗹𝗼𝗰𝘀 ✧ (𝘀𝗲𝗹𝗳) 𝗹𝗶𝗯𝗲𝗿𝗮𝘁𝗶𝗼𝗻 ✧ 𝗹𝗶𝗳𝗲𝘀𝘁𝘆𝗹𝗲 ✧
🪴 | 71 ʟᴏᴄs est. on 07.20.23
📍 | ʜᴏᴜsᴛᴏɴ, ᴛx"
$ws.Cells.Item(42, 2).Value = 2979
$ws.Cells.Item(42, 3).Value = 238
$ws.Cells.Item(42, 4).Value = 172
$ws.Cells.Item(42, 5).Value = "https://scontent-hou1-1.xx.fbcdn.net/v/t51.2885-15/481266977_997353345602937_1719041919639027270_n.jpg?_nc_cat=106&ccb=1-7&_nc_sid=7d201b&_nc_ohc=98Z30Y-Fn84Q7kNvwGigcY3&_nc_oc=AdnlmqP0Owtm9g4fPHIQtb0Rs-3ny1ApT5-BOO0J8fyh8uZsEbvN01ZdKBP2dzXKBWBYiGD6BG8VmcDe86MNKr2X&_nc_zt=23&_nc_ht=scontent-hou1-1.xx&edm=AL-3X8kEAAAA&oh=00_AfI74XBNz4sZmfpTCCRLRUo6v8FSnNHdPQSIPqEIpumutQ&oe=6826806E"
$ws.Cells.Item(42, 6).Value = "'17841461458191255"
$ws.Cells.Item(42, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(42, 7).Value = 45788.39804859953
$ws.Cells.Item(42, 8).Value = "'2025-05-11"
$ws.Cells.Item(42, 9).Value = 2025
$ws.Cells.Item(42, 10).Value = "May"
$ws.Cells.Item(42, 11).Value = 11
$ws.Cells.Item(42, 12).Value = "09:33:11"

# --- Row 43 ---
$ws.Cells.Item(43, 1).Value = "ʜᴇʟᴘɪɴɢ ʏᴏᴜ ᴍᴏᴠᴇ ғʀᴏᴍ ʜᴇsɪᴛᴀᴛɪᴏɴ ᴛᴏ ᴄʀᴇᴀᴛɪᴏɴ
✨ | ✧ 𝗹𝗼𝗰𝘀 ✧ (𝘀𝗲𝗹𝗳) 𝗹𝗶𝗯𝗲𝗿𝗮𝘁𝗶𝗼𝗻 ✧ 𝗹𝗶𝗳𝗲𝘀𝘁𝘆𝗹𝗲 ✧
🪴 | 71 ʟᴏᴄs est. on 07.20.23
📍 | ʜᴏᴜsᴛᴏɴ, ᴛx"
$ws.Cells.Item(43, 2).Value = 2983
$ws.Cells.Item(43, 3).Value = 238
$ws.Cells.Item(43, 4).Value = 172
$ws.Cells.Item(43, 5).Value = "https://scontent-hou1-1.xx.fbcdn.net/v/t51.2885-15/481266977_997353345602937_1719041919639027270_n.jpg?_nc_cat=106&ccb=1-7&_nc_sid=7d201b&_nc_ohc=98Z30Y-Fn84Q7kNvwFG_Vub&_nc_oc=AdlBhQs5yoAx2HX9ehwQEDxKusa4_dxclDLJsZm2zFF7LvHd2ATdglLLxKuRHjvPxDybk4J6fG9sXmCu5a9P_MPF&_nc_zt=23&_nc_ht=scontent-hou1-1.xx&edm=AL-3X8kEAAAA&oh=00_AfIFkVgVIEsFPA10mlaTSX0D9GLV76yK-VCA3k6ArQ6NAw&oe=6827D1EE"
$ws.Cells.Item(43, 6).Value = "'17841461458191255"
$ws.Cells.Item(43, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(43, 7).Value = 45789.42137029266
$ws.Cells.Item(43, 8).Value = "'2025-05-12"
$ws.Cells.Item(43, 9).Value = 2025
$ws.Cells.Item(43, 10).Value = "May"
$ws.Cells.Item(43, 11).Value = 12
$ws.Cells.Item(43, 12).Value = "10:06:46"

